$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells
$ws.Range("H1").Value = "Tempo Heuristica"
$ws.Range("I1").Value = "Tempo Total"

# Update G column (Tempo) values and add new H (Tempo Heuristica) / I (Tempo Total) values
$ws.Range("G2").Value = 0.0569303035736084
$ws.Range("H2").Value = 0.009218215942382812
$ws.Range("I2").Value = 0.06614851951599121

$ws.Range("G3").Value = 0.05357766151428223
$ws.Range("H3").Value = 0.01400089263916016
$ws.Range("I3").Value = 0.06757855415344238

$ws.Range("G4").Value = 0.05604004859924316
$ws.Range("H4").Value = 0.009263277053833008
$ws.Range("I4").Value = 0.06530332565307617
